# Sustainable Forms presentation.pptx - slide 6 edits
#
# 1. Title shape ("Lägesrapport produktutveckling") -> "Lägesrapport utveckling"
# 2. Sub-title textbox ("Var är vi nu?") -> cleared (text removed)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- 1. Update the title text ---
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Lägesrapport utveckling"

# --- 2. Clear the "Var är vi nu?" textbox ---
$subtitle = $s.Shapes.Item(3)
$subtitle.TextFrame.TextRange.Delete()
